$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("G1").Value = "Pathway"
$ws.Range("H1").Value = "Hemi"

# ventral_right: rows 2-13
$ws.Range("G2:G13").Value = "ventral"
$ws.Range("H2:H13").Value = "right"

# ventral_left: rows 14-25
$ws.Range("G14:G25").Value = "ventral"
$ws.Range("H14:H25").Value = "left"

# dorsal_right: rows 26-37
$ws.Range("G26:G37").Value = "dorsal"
$ws.Range("H26:H37").Value = "right"

# dorsal_left: rows 38-49
$ws.Range("G38:G49").Value = "dorsal"
$ws.Range("H38:H49").Value = "left"
